# Fruta / hortaliza, semanal
# Insert a new data row at row 85 (pushing the existing rows 85-88 down to
# 86-89) and populate it with a new weekly price record for
# Femacal de La Calera / Tuna.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 85..88 down to 86..89, leaving a fresh blank row 85.
$ws.Rows.Item(85).Insert()

$ws.Range("A85").Value = 3
$ws.Range("B85").Value = 'Femacal de La Calera'
$ws.Range("C85").Value = 'Coquimbo'
$ws.Range("D85").Value = 44585
$ws.Range("E85").Value = 5
$ws.Range("F85").Value = 'Fruta'
$ws.Range("G85").Value = 100107
$ws.Range("H85").Value = 'Otros'
$ws.Range("I85").Value = 100107011
$ws.Range("J85").Value = 'Tuna'
$ws.Range("K85").Value = 'Sin especificar'
$ws.Range("L85").Value = 'Primera'
$ws.Range("M85").Value = 50
$ws.Range("N85").Value = 20000
$ws.Range("O85").Value = 20000
$ws.Range("P85").Value = 20000
$ws.Range("Q85").Value = '$/caja 20 kilos'
$ws.Range("R85").Value = 'Provincia de Limarí'
$ws.Range("S85").Value = 1000
$ws.Range("T85").Value = 20
